$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.736.37"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "3.731.38"
$ws.Range("E3").Value = "  +18.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.39"
$ws.Range("E5").Value = "  +5.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.01"
$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").Value = "3.724.37"
$ws.Range("E7").Value = "  +18.30%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  +3.85%  "

$ws.Range("E10").Value = "  +9.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  +6.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.99"
$ws.Range("E13").Value = "  +10.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  +5.25%  "

$ws.Range("D15").Value = "4.354.83"
$ws.Range("E15").Value = "  +18.54%  "

$ws.Range("D16").Value = "3.731.18"
$ws.Range("E16").Value = "  +18.61%  "

$ws.Range("D17").Value = "69.791.06"
$ws.Range("E17").Value = "  +1.84%  "

$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  +6.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "517.41"
$ws.Range("E20").Value = "  +5.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.70"
$ws.Range("E21").Value = "  +1.40%  "

$ws.Range("E22").Value = "  +19.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.731"
$ws.Range("E23").Value = "  +4.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.51"
$ws.Range("E24").Value = "  +5.35%  "

$ws.Range("E25").Value = "  +5.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.60"
$ws.Range("E26").Value = "  +4.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  +3.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  +32.76%  "

$ws.Range("E30").Value = "  +6.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.85"
$ws.Range("E31").Value = "  +7.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  -3.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.45"
$ws.Range("E33").Value = "  +11.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.22"
$ws.Range("E36").Value = "  +7.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("E37").Value = "  +8.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.341"
$ws.Range("E38").Value = "  +4.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.19"
$ws.Range("E39").Value = "  +6.75%  "

$ws.Range("E40").Value = "  +5.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.30"
$ws.Range("E41").Value = "  +4.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.78"
$ws.Range("E42").Value = "  -7.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("E43").Value = "  +4.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "424.06"
$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("D45").Value = "3.058.79"
$ws.Range("E45").Value = "  +8.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("E47").Value = "  +4.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.94"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  +6.04%  "

$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.59"
$ws.Range("E51").Value = "  +0.66%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"